# Bug Fixes in local attachments with ignore file
#
# Appends the latest UI education registration entries to the
# "AMSIN" and "AMS" history sheets:
#   - AMSIN: normalizes the previously "latest" row (57) to the
#     sheet's standard formatting and appends the new latest row (58).
#   - AMS: appends two new rows (36, 37); row 36 gets the standard
#     formatting (it is no longer the latest row) while row 37 (the
#     new latest row) is left with default/no explicit formatting -
#     matching how every sheet's freshly-appended row starts out.

$wb = $excel.ActiveWorkbook

# Date/time number format used throughout these sheets for column B.
$dtFormat = "YYYY-MM-DD HH:MM:SS"

function Set-TextCell($cell, $text) {
    # Force literal text storage so date-looking strings (e.g. "2023-07-31")
    # aren't reinterpreted by Excel as date serials.
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

function Set-DateTimeCell($cell, $serial) {
    $cell.NumberFormat = $dtFormat
    $cell.Value = $serial
}

# ---------------------------------------------------------------------
# Sheet "AMSIN": normalize row 57, append row 58
# (columns default to the sheet's standard style, so clearing a cell
# before rewriting its value restores that standard look)
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

# Row 57 already holds its data; only its formatting needs to catch up
# to match the rest of the sheet.
$c = $wsAmsin.Cells.Item(57, 1)
$c.Clear()
Set-TextCell $c "2023-06-12"
$wsAmsin.Cells.Item(57, 2).Value = 45089.6148461574
foreach ($col in 3,4,5,6,7) {
    $wsAmsin.Cells.Item(57, $col).Clear()
}
$wsAmsin.Cells.Item(57, 3).Value = "178ddyedu"
$wsAmsin.Cells.Item(57, 4).Value = 59
$wsAmsin.Cells.Item(57, 5).Value = 59
$wsAmsin.Cells.Item(57, 6).Value = 0
$wsAmsin.Cells.Item(57, 7).Value = 1.07

# Row 58: brand new latest entry.
Set-TextCell $wsAmsin.Cells.Item(58, 1) "2023-07-31"
Set-DateTimeCell $wsAmsin.Cells.Item(58, 2) 45138.42191619213
$wsAmsin.Cells.Item(58, 3).Value = "180educ"
$wsAmsin.Cells.Item(58, 4).Value = 59
$wsAmsin.Cells.Item(58, 5).Value = 59
$wsAmsin.Cells.Item(58, 6).Value = 0
$wsAmsin.Cells.Item(58, 7).Value = 0.85

# ---------------------------------------------------------------------
# Sheet "AMS": append rows 36 and 37
# (this sheet has no column-level default style, so the standard look
# for row 36 is pulled explicitly from the row above it)
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Row 36: gets the sheet's standard formatting (row 35 is the template).
Set-TextCell $wsAms.Cells.Item(36, 1) "2023-08-01"
Set-DateTimeCell $wsAms.Cells.Item(36, 2) 45139.5405990625
$wsAms.Cells.Item(36, 3).Value = "180educa"
$wsAms.Cells.Item(36, 4).Value = 59
$wsAms.Cells.Item(36, 5).Value = 59
$wsAms.Cells.Item(36, 6).Value = 0
$wsAms.Cells.Item(36, 7).Value = 0.89
foreach ($col in 3,4,5,6,7) {
    $dst = $wsAms.Cells.Item(36, $col)
    $src = $wsAms.Cells.Item(35, $col)
    $dst.Style = $src.Style
}

# Row 37: newest row, left unformatted like every sheet's latest row.
Set-TextCell $wsAms.Cells.Item(37, 1) "2023-08-01"
Set-DateTimeCell $wsAms.Cells.Item(37, 2) 45139.86483763212
$wsAms.Cells.Item(37, 3).Value = "180liveeuc"
$wsAms.Cells.Item(37, 4).Value = 59
$wsAms.Cells.Item(37, 5).Value = 59
$wsAms.Cells.Item(37, 6).Value = 0
$wsAms.Cells.Item(37, 7).Value = 0.86
